# Update column K (IRE) values for rows 5-8 with newly scraped batting data.
# New value scraped for K5 (4359) shifts the previously existing values down
# one row: old K5 -> K6, old K6 -> K7, old K7 -> K8 (which was previously empty).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("K5:K8")
# Ensure the cells keep a text format so purely-numeric strings are not
# auto-converted to the Number type (matches original inline-string cells).
$rng.NumberFormat = "@"

$ws.Range("K5").Value = "4359"
$ws.Range("K6").Value = "3418"
$ws.Range("K7").Value = "6148"
$ws.Range("K8").Value = "4586"

# Reset the style back to Normal/default so no residual styling (beyond the
# text format) is applied to these cells, matching the rest of the sheet.
$rng.Style = "Normal"
